$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new "% GO spectra - cellular component" values (I3:L3)
$ws.Range("I3").Value = 0.0766632894292468
$ws.Range("J3").Value = 0.0648429584599797
$ws.Range("K3").Value = 0.180344478216818
$ws.Range("L3").Value = 0.0536980749746707

# Row 4
$ws.Range("I4").Value = 0.123188405797101
$ws.Range("J4").Value = 0.108695652173913
$ws.Range("K4").Value = 0.137681159420289
$ws.Range("L4").Value = 0.123188405797101

# Row 5
$ws.Range("I5").Value = 0.07803993
$ws.Range("J5").Value = 0.05444646
$ws.Range("K5").Value = 0.08529946
$ws.Range("L5").Value = 0.04900182

# Row 6
$ws.Range("I6").Value = 0.13357731
$ws.Range("J6").Value = 0.08691674
$ws.Range("K6").Value = 0.11619396
$ws.Range("L6").Value = 0.09240622

# Update the saved selection to N3 (was H7)
$ws.Range("N3").Select()

# Touch page setup so an explicit <pageSetup> element is written (portrait orientation)
$ws.PageSetup.Orientation = 1
